$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A28").Value = 46005
$ws.Range("B28").Value = "四方坪站"
$ws.Range("C28").Value = 8725.93
$ws.Range("D28").Value = 7279.61
$ws.Range("E28").Value = 2907.75
$ws.Range("F28").Value = 354

$ws.Range("A29").Value = 46005
$ws.Range("B29").Value = "高岭站"
$ws.Range("C29").Value = 5263.11
$ws.Range("D29").Value = 4406.77
$ws.Range("E29").Value = 1353.99
$ws.Range("F29").Value = 181

$ws.Range("H27").Select()
